$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "/home/daniel/Spike Data/Matlab files/Exp 27 unit 1 data.mat"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 157400
$ws.Range("F9").Value = 167200

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "/home/daniel/Spike Data/Matlab files/31 Slow ramp.mat"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 4027
$ws.Range("F10").Value = 11550

$ws.Range("C16").Select()
